# Update p_ttes_c_charge_discharge (column H, row 2) from 5.5 to 5
# on every year-sheet except 2025 (which already holds the value 5).
$wb = $excel.ActiveWorkbook

$sheetNames = @("2030", "2035", "2040", "2045", "2050")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("H2").Value = 5
}
